# Implementacion final de notificacion por correo (Enviado, Observado)
#
# Adds a new "notification status" column (E) to the AlumnoCursoTesis2022
# sheet, defaulting every existing student row (1-90) to "null", and moves
# the active selection to G6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row in column A (the sheet currently holds 90 student
# rows, A1:D90) so the new column is populated for every existing row.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 1; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 5).Value = "null"
}

# Match the author's final cursor position after the edit.
[void]$ws.Range("G6").Select()
